# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D) for the acfd965a... file
# row (row 4) on both the zh-cn and de-de localization-status sheets,
# reflecting a newer handoff that just completed.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-03-04 08:08:24"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-03-04 08:08:33"
